# Amend corrected label annotations
# Normalize the "labels" column (F) values to lowercase for the rows
# identified in the diff. A couple of multi-label rows also have their
# "||"-joined segments reordered to match the corrected annotation order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5   = "18_hazards_to_humans_and_domestic_animals"
    6   = "ppe"
    9   = "off target movement || application instructions || env warning - species || env warning - water"
    13  = "application instructions"
    14  = "135_product_information"
    15  = "application instructions"
    67  = "use restrictions"
    75  = "application instructions"
    77  = "off target movement"
    78  = "172_sensitive_areas"
    79  = "off target movement"
    80  = "off target movement"
    81  = "application instructions || off target movement"
    82  = "off target movement"
    83  = "off target movement"
    84  = "off target movement"
    85  = "application instructions"
    86  = "application instructions"
    87  = "mixing"
    88  = "mixing"
    89  = "mixing"
    90  = "mixing"
    91  = "mixing"
    92  = "safety procedures"
    93  = "use restrictions"
    94  = "use restrictions"
    97  = "application instructions"
    98  = "application instructions"
    99  = "use restrictions"
    100 = "mixing"
    101 = "mixing"
    104 = "mixing"
    105 = "mixing"
    108 = "mixing"
    109 = "mixing"
    125 = "mixing"
    127 = "mixing"
    131 = "mixing"
    135 = "mixing"
    136 = "mixing"
    137 = "mixing"
    138 = "application instructions"
    139 = "application instructions"
    140 = "use restrictions"
    141 = "mixing"
    143 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
